$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8952653408050537
$ws.Range("B1").Value = 1.327225923538208
$ws.Range("C1").Value = 2.749033451080322
$ws.Range("D1").Value = 3.612663984298706
$ws.Range("E1").Value = 1.769871115684509
